# Auto-generated Excel COM-interop script applying the Spriggan_Profits diff
# Updates market-price / profit columns (H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 9263012
$ws.Range("I40").Value = 3003.7
$ws.Range("J40").Value = 55563056
$ws.Range("K40").Value = 3003.7
$ws.Range("L40").Value = 55563056
$ws.Range("M40").Value = -2828.7
$ws.Range("N40").Value = -55563406

$ws.Range("H69").Value = 16992.895
$ws.Range("I69").Value = 7717.25
$ws.Range("J69").Value = 19466.4
$ws.Range("K69").Value = 23151.75
$ws.Range("L69").Value = 58399.2
$ws.Range("M69").Value = -22277.75
$ws.Range("N69").Value = -60147.2

$ws.Range("H72").Value = 16992.895
$ws.Range("I72").Value = 7717.25
$ws.Range("J72").Value = 19466.4
$ws.Range("K72").Value = 69455.25
$ws.Range("L72").Value = 175197.6
$ws.Range("M72").Value = -65087.25
$ws.Range("N72").Value = -183933.6

$ws.Range("H100").Value = 2584.077
$ws.Range("I100").Value = 2274.25
$ws.Range("J100").Value = 3079.8
$ws.Range("K100").Value = 2274.25
$ws.Range("L100").Value = 3079.8
$ws.Range("M100").Value = -1733.25
$ws.Range("N100").Value = -4161.8

$ws.Range("H113").Value = 2658.8
$ws.Range("I113").Value = 2698.5
$ws.Range("K113").Value = 2698.5
$ws.Range("M113").Value = 555.5

$ws.Range("H138").Value = 6674.4614
$ws.Range("I138").Value = 5428.6665
$ws.Range("J138").Value = 7048.2
$ws.Range("K138").Value = 16285.9995
$ws.Range("L138").Value = 21144.6
$ws.Range("M138").Value = -11145.9995
$ws.Range("N138").Value = -31424.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2131.9253
$ws.Range("I32").Value = 981.875
$ws.Range("K32").Value = 981.875
$ws.Range("M32").Value = -694.875

$ws.Range("H97").Value = 831.86664
$ws.Range("I97").Value = 888.8182
$ws.Range("K97").Value = 888.8182
$ws.Range("M97").Value = -392.8182

$ws.Range("H132").Value = 6253296
$ws.Range("I132").Value = 8336637
$ws.Range("K132").Value = 25009911
$ws.Range("M132").Value = -25007381

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 39115.25
$ws.Range("I107").Value = 3531.0476
$ws.Range("J107").Value = 145867.86
$ws.Range("K107").Value = 3531.0476
$ws.Range("L107").Value = 145867.86
$ws.Range("M107").Value = -1611.0476
$ws.Range("N107").Value = -149707.86

$ws.Range("H134").Value = 56668280
$ws.Range("I134").Value = 72858344
$ws.Range("J134").Value = 3049.5
$ws.Range("K134").Value = 218575032
$ws.Range("L134").Value = 9148.5
$ws.Range("M134").Value = -218572497
$ws.Range("N134").Value = -14218.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3530.804
$ws.Range("I31").Value = 5403.1665
$ws.Range("J31").Value = 2509.5151
$ws.Range("K31").Value = 5403.1665
$ws.Range("L31").Value = 2509.5151
$ws.Range("M31").Value = -5108.1665
$ws.Range("N31").Value = -3099.5151

$ws.Range("H34").Value = 3530.804
$ws.Range("I34").Value = 5403.1665
$ws.Range("J34").Value = 2509.5151
$ws.Range("K34").Value = 5403.1665
$ws.Range("L34").Value = 2509.5151
$ws.Range("M34").Value = -5201.1665
$ws.Range("N34").Value = -2913.5151

$ws.Range("H86").Value = 11186.056
$ws.Range("I86").Value = 9194.299999999999
$ws.Range("J86").Value = 13675.75
$ws.Range("K86").Value = 9194.299999999999
$ws.Range("L86").Value = 13675.75
$ws.Range("M86").Value = -8071.299999999999
$ws.Range("N86").Value = -15921.75

$ws.Range("H89").Value = 11186.056
$ws.Range("I89").Value = 9194.299999999999
$ws.Range("J89").Value = 13675.75
$ws.Range("K89").Value = 45971.5
$ws.Range("L89").Value = 68378.75
$ws.Range("M89").Value = -40355.5
$ws.Range("N89").Value = -79610.75

$ws.Range("H94").Value = 2366.1667
$ws.Range("I94").Value = 2466.3333
$ws.Range("K94").Value = 2466.3333
$ws.Range("M94").Value = -2015.3333

$ws.Range("H107").Value = 505024.03
$ws.Range("I107").Value = 529370.4399999999
$ws.Range("J107").Value = 334599.34
$ws.Range("K107").Value = 529370.4399999999
$ws.Range("L107").Value = 334599.34
$ws.Range("M107").Value = -527450.4399999999
$ws.Range("N107").Value = -338439.34

$ws.Range("H132").Value = 18185182
$ws.Range("I132").Value = 20836544
$ws.Range("K132").Value = 62509632
$ws.Range("M132").Value = -62507102

$ws.Range("H134").Value = 6100708.5
$ws.Range("I134").Value = 6252976
$ws.Range("J134").Value = 9999
$ws.Range("K134").Value = 18758928
$ws.Range("L134").Value = 29997
$ws.Range("M134").Value = -18756393
$ws.Range("N134").Value = -35067

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 845
$ws.Range("I22").Value = 700
$ws.Range("J22").Value = 990
$ws.Range("K22").Value = 2100
$ws.Range("L22").Value = 2970
$ws.Range("M22").Value = -1931
$ws.Range("N22").Value = -3308

$ws.Range("H27").Value = 845
$ws.Range("I27").Value = 700
$ws.Range("J27").Value = 990
$ws.Range("K27").Value = 2100
$ws.Range("L27").Value = 2970
$ws.Range("M27").Value = -1998
$ws.Range("N27").Value = -3174

$ws.Range("H68").Value = 7142.778
$ws.Range("J68").Value = 7206.4917
$ws.Range("L68").Value = 21619.4751
$ws.Range("N68").Value = -23241.4751

$ws.Range("H71").Value = 7142.778
$ws.Range("J71").Value = 7206.4917
$ws.Range("L71").Value = 64858.4253
$ws.Range("N71").Value = -72970.4253

$ws.Range("H122").Value = 793.35297
$ws.Range("I122").Value = 695.8
$ws.Range("J122").Value = 932.7143
$ws.Range("K122").Value = 6262.2
$ws.Range("L122").Value = 8394.4287
$ws.Range("M122").Value = -3812.2
$ws.Range("N122").Value = -13294.4287

$ws.Range("H129").Value = 1662.5385
$ws.Range("I129").Value = 1004.7895
$ws.Range("J129").Value = 3447.8572
$ws.Range("K129").Value = 3014.3685
$ws.Range("L129").Value = 10343.5716
$ws.Range("M129").Value = 1985.6315
$ws.Range("N129").Value = -20343.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 34496.332
$ws.Range("J32").Value = 34496.332
$ws.Range("L32").Value = 34496.332
$ws.Range("N32").Value = -35088.332

$ws.Range("H107").Value = 2067.5
$ws.Range("I107").Value = 1457.0769
$ws.Range("J107").Value = 10003
$ws.Range("K107").Value = 1457.0769
$ws.Range("L107").Value = 10003
$ws.Range("M107").Value = 462.9231
$ws.Range("N107").Value = -13843

$ws.Range("H122").Value = 119090.91
$ws.Range("I122").Value = 174571.42
$ws.Range("K122").Value = 523714.26
$ws.Range("M122").Value = -521264.26

$ws.Range("H126").Value = 3100.5
$ws.Range("I126").Value = 3172.9333
$ws.Range("J126").Value = 2014
$ws.Range("K126").Value = 9518.7999
$ws.Range("L126").Value = 6042
$ws.Range("M126").Value = -7048.7999
$ws.Range("N126").Value = -10982

$ws.Range("H132").Value = 7356452.5
$ws.Range("I132").Value = 9618454
$ws.Range("J132").Value = 4947.25
$ws.Range("K132").Value = 28855362
$ws.Range("L132").Value = 14841.75
$ws.Range("M132").Value = -28852832
$ws.Range("N132").Value = -19901.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3782.5
$ws.Range("I7").Value = 3839
$ws.Range("J7").Value = 3500
$ws.Range("K7").Value = 3839
$ws.Range("L7").Value = 3500
$ws.Range("M7").Value = -3727
$ws.Range("N7").Value = -3724

$ws.Range("H16").Value = 3898.2856
$ws.Range("I16").Value = 3486.3333
$ws.Range("J16").Value = 4639.8
$ws.Range("K16").Value = 3486.3333
$ws.Range("L16").Value = 4639.8
$ws.Range("M16").Value = -3316.3333
$ws.Range("N16").Value = -4979.8

$ws.Range("H55").Value = 418.85
$ws.Range("I55").Value = 119.36364
$ws.Range("J55").Value = 784.8889
$ws.Range("K55").Value = 119.36364
$ws.Range("L55").Value = 784.8889
$ws.Range("M55").Value = 53.63636
$ws.Range("N55").Value = -1130.8889

$ws.Range("H61").Value = 3374.3572
$ws.Range("I61").Value = 3380.4443
$ws.Range("K61").Value = 3380.4443
$ws.Range("M61").Value = -3178.4443

$ws.Range("H68").Value = 1394710.6
$ws.Range("I68").Value = 2928809.5
$ws.Range("J68").Value = 14021.5
$ws.Range("K68").Value = 2928809.5
$ws.Range("L68").Value = 14021.5
$ws.Range("M68").Value = -2928060.5
$ws.Range("N68").Value = -15519.5

$ws.Range("H71").Value = 1394710.6
$ws.Range("I71").Value = 2928809.5
$ws.Range("J71").Value = 14021.5
$ws.Range("K71").Value = 14644047.5
$ws.Range("L71").Value = 70107.5
$ws.Range("M71").Value = -14640303.5
$ws.Range("N71").Value = -77595.5

$ws.Range("H93").Value = 2125.3
$ws.Range("I93").Value = 1794.4706
$ws.Range("J93").Value = 4000
$ws.Range("K93").Value = 1794.4706
$ws.Range("L93").Value = 4000
$ws.Range("M93").Value = -546.4706000000001
$ws.Range("N93").Value = -6496

$ws.Range("H113").Value = 3374.3572
$ws.Range("I113").Value = 3380.4443
$ws.Range("K113").Value = 3380.4443
$ws.Range("M113").Value = -1210.4443

$ws.Range("H126").Value = 3782.5
$ws.Range("I126").Value = 3839
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 11517
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -9047
$ws.Range("N126").Value = -15440

$ws.Range("H132").Value = 7815944.5
$ws.Range("I132").Value = 10420458
$ws.Range("J132").Value = 2404.625
$ws.Range("K132").Value = 31261374
$ws.Range("L132").Value = 7213.875
$ws.Range("M132").Value = -31258844
$ws.Range("N132").Value = -12273.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 70024.734
$ws.Range("I81").Value = 80036.84
$ws.Range("J81").Value = 4946
$ws.Range("K81").Value = 160073.68
$ws.Range("L81").Value = 9892
$ws.Range("M81").Value = -159012.68
$ws.Range("N81").Value = -12014

$ws.Range("H84").Value = 70024.734
$ws.Range("I84").Value = 80036.84
$ws.Range("J84").Value = 4946
$ws.Range("K84").Value = 800368.3999999999
$ws.Range("L84").Value = 49460
$ws.Range("M84").Value = -795064.3999999999
$ws.Range("N84").Value = -60068

$ws.Range("H113").Value = 1292
$ws.Range("I113").Value = 979.1111
$ws.Range("J113").Value = 2230.6667
$ws.Range("K113").Value = 2937.3333
$ws.Range("L113").Value = 6692.000100000001
$ws.Range("M113").Value = -767.3332999999998
$ws.Range("N113").Value = -11032.0001

$ws.Range("H122").Value = 1529.75
$ws.Range("I122").Value = 1045.65
$ws.Range("J122").Value = 2740
$ws.Range("K122").Value = 3136.95
$ws.Range("L122").Value = 8220
$ws.Range("M122").Value = -686.9500000000003
$ws.Range("N122").Value = -13120

$ws.Range("H129").Value = 94999.5
$ws.Range("J129").Value = 94999.5
$ws.Range("L129").Value = 94999.5
$ws.Range("N129").Value = -104999.5

